$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 584; this pushes the
# existing rows 584:643 down to 585:644 (and grows the used range /
# dimension from R643 to R644 automatically).
$ws.Rows("584:584").Insert()

# Populate the newly inserted row 584 with this week's record. Static
# "template" columns (market, region, category, quality, unit, origin,
# classification, etc.) mirror the row that used to sit at 584 before the
# shift; only the date / volume / price columns carry genuinely new data.
$ws.Cells.Item(584, 1).Value2 = 6
$ws.Cells.Item(584, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(584, 3).Value2 = "Metropolitana"
$ws.Cells.Item(584, 4).Value2 = 44918
$ws.Cells.Item(584, 5).Value2 = 13
$ws.Cells.Item(584, 6).Value2 = 100112044
$ws.Cells.Item(584, 7).Value2 = "Perejil"
$ws.Cells.Item(584, 8).Value2 = "Sin especificar"
$ws.Cells.Item(584, 9).Value2 = "Primera"
$ws.Cells.Item(584, 10).Value2 = 340
$ws.Cells.Item(584, 11).Value2 = 10000
$ws.Cells.Item(584, 12).Value2 = 11000
$ws.Cells.Item(584, 13).Value2 = 10441
$ws.Cells.Item(584, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(584, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(584, 16).Value2 = 3480
$ws.Cells.Item(584, 17).Value2 = 3
$ws.Cells.Item(584, 18).Value2 = "Hortaliza"

# Keep the date formatting consistent with the rest of column D.
$ws.Cells.Item(584, 4).NumberFormat = $ws.Cells.Item(585, 4).NumberFormat
